$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.475.64"
$ws.Range("E2").Value = "  -0.37%  "

$ws.Range("D3").Value = "2.661.93"
$ws.Range("E3").Value = "  +0.97%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.45%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.05%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  -0.84%  "

$ws.Range("D9").Value = "2.660.28"
$ws.Range("E9").Value = "  +0.94%  "

$ws.Range("E10").Value = "  -4.33%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.169"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.85%  "

$ws.Range("E12").Value = "  +0.03%  "

$ws.Range("E13").Value = "  -0.63%  "

$ws.Range("D14").Value = "3.149.21"
$ws.Range("E14").Value = "  +1.33%  "

$ws.Range("D15").Value = "72.401.01"
$ws.Range("E15").Value = "  -0.21%  "

$ws.Range("E16").Value = "  -3.19%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.22"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.27%  "

$ws.Range("D18").Value = "2.653.34"
$ws.Range("E18").Value = "  +0.98%  "

$ws.Range("E19").Value = "  +5.31%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.85%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "370.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.91%  "

$ws.Range("E22").Value = "  +0.32%  "

$ws.Range("E23").Value = "  +2.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.04%  "

$ws.Range("E25").Value = "  -0.12%  "

$ws.Range("E26").Value = "  -2.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.75%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.04%  "

$ws.Range("D30").Value = "0.0₃0969"
$ws.Range("E30").Value = "  +0.79%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.30%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "495.77"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.52%  "

$ws.Range("E33").Value = "  -2.31%  "

$ws.Range("E34").Value = "  -0.34%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.06%  "

$ws.Range("E36").Value = "  -0.02%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.50"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.40%  "

$ws.Range("E38").Value = "  +0.54%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.93"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.96%  "

$ws.Range("E40").Value = "  -2.06%  "

$ws.Range("E41").Value = "  -5.38%  "

$ws.Range("E42").Value = "  -0.07%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.99"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.18%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.58"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.48%  "

$ws.Range("E45").Value = "  -0.43%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "156.63"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.56%  "

$ws.Range("E47").Value = "  -0.42%  "

$ws.Range("E48").Value = "  +2.34%  "

$ws.Range("E49").Value = "  +0.92%  "

$ws.Range("E50").Value = "  +1.60%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0758"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.25%  "
